# Edit script: apply diff changes to Technical Solution.docx
$d = $word.ActiveDocument

# --- 1. Paragraph 2 ("Completeness of Solution - ...") : remove the _GoBack bookmark ---
$p2 = $d.Paragraphs.Item(2)
$xmlPara2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3A592867" w14:textId="59092750" w:rsidR="000A0945" w:rsidRPr="003937C8" w:rsidRDefault="004368F9" w:rsidP="00C92711" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Completeness of Solution - </w:t></w:r><w:r w:rsidR="00D12D3D" w:rsidRPr="003937C8"><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t>A system that meets almost all of the requirements of a solution/an investigation (ignoring any requirements that go beyond the demands of A-level).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xmlPara2)

# --- 2. Paragraph 5 ("Above average performance...") : split "Group A" with proofErr marks ---
$p5 = $d.Paragraphs.Item(5)
$xmlPara5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="7D8E3005" w14:textId="2B616539" w:rsidR="003937C8" w:rsidRPr="00643A19" w:rsidRDefault="003937C8" w:rsidP="003937C8" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Techniques Used Additional Information - </w:t></w:r><w:r w:rsidRPr="00643A19"><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Above average performance: Group </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> equivalent algorithms and model programmed more than well to excellent; all or almost all excellent coding style characteristics; more than to highly effective solution.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($xmlPara5)

# --- 3. Paragraph 6 ("Average performance...") : split "Group A" with proofErr marks ---
$p6 = $d.Paragraphs.Item(6)
$xmlPara6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="52D51B30" w14:textId="77777777" w:rsidR="003937C8" w:rsidRPr="00643A19" w:rsidRDefault="003937C8" w:rsidP="003937C8" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00643A19"><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Average performance: Group </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> equivalent algorithms and/or model programmed well; majority of excellent coding style characteristics; an effective solution.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p6.Range.InsertXML($xmlPara6)

# --- 4. Paragraph 7 ("Below average performance...") : split "Group A" with proofErr marks,
#        and move the _GoBack bookmark to the end of this paragraph ---
$p7 = $d.Paragraphs.Item(7)
$xmlPara7 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="0D0C2732" w14:textId="6FF1D410" w:rsidR="003937C8" w:rsidRDefault="003937C8" w:rsidP="003937C8" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00643A19"><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Below average performance: Group </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="3B3B3A"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> equivalent algorithms and/or model programmed just adequately to fully adequate; some excellent coding style characteristics; less than effective to fairly effective solution.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p7.Range.InsertXML($xmlPara7)

# --- 5. Replace the trailing empty paragraph with the new content block ---
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$xmlTail = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>??? Ask Steven</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Techniques Used to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Create</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>mmltomidi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Techniques Used to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Create</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>catmidi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pLast.Range.InsertXML($xmlTail)

# InsertXML with multiple <w:p> elements inserts them *before* the target paragraph,
# leaving the original (now-redundant) empty paragraph dangling at the end of the body.
# Remove that leftover empty paragraph by merging its paragraph mark into the
# paragraph immediately before it.
$newCount = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($newCount - 1)
$veryLast = $d.Paragraphs.Item($newCount)
$cleanupRange = $d.Range($secondLast.Range.End - 1, $veryLast.Range.End)
$cleanupRange.Delete()
